$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - set value then copy formatting from E1 (bold/centered/bordered header style)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2:F98 - time_taken metadata timestamps, default (unstyled) format
$ws.Range("F2").Value = "2021-10-05 13:40:02.997785"
$ws.Range("F3").Value = "2021-10-05 13:40:02.997796"
$ws.Range("F4").Value = "2021-10-05 13:40:02.997800"
$ws.Range("F5").Value = "2021-10-05 13:40:02.997802"
$ws.Range("F6").Value = "2021-10-05 13:40:02.997805"
$ws.Range("F7").Value = "2021-10-05 13:40:02.997808"
$ws.Range("F8").Value = "2021-10-05 13:40:02.997810"
$ws.Range("F9").Value = "2021-10-05 13:40:02.997813"
$ws.Range("F10").Value = "2021-10-05 13:40:02.997816"
$ws.Range("F11").Value = "2021-10-05 13:40:02.997818"
$ws.Range("F12").Value = "2021-10-05 13:40:02.997821"
$ws.Range("F13").Value = "2021-10-05 13:40:02.997824"
$ws.Range("F14").Value = "2021-10-05 13:40:02.997826"
$ws.Range("F15").Value = "2021-10-05 13:40:02.997829"
$ws.Range("F16").Value = "2021-10-05 13:40:02.997831"
$ws.Range("F17").Value = "2021-10-05 13:40:02.997833"
$ws.Range("F18").Value = "2021-10-05 13:40:02.997836"
$ws.Range("F19").Value = "2021-10-05 13:40:02.997839"
$ws.Range("F20").Value = "2021-10-05 13:40:02.997841"
$ws.Range("F21").Value = "2021-10-05 13:40:02.997844"
$ws.Range("F22").Value = "2021-10-05 13:40:02.997846"
$ws.Range("F23").Value = "2021-10-05 13:40:02.997849"
$ws.Range("F24").Value = "2021-10-05 13:40:02.997851"
$ws.Range("F25").Value = "2021-10-05 13:40:02.997854"
$ws.Range("F26").Value = "2021-10-05 13:40:02.997857"
$ws.Range("F27").Value = "2021-10-05 13:40:02.997859"
$ws.Range("F28").Value = "2021-10-05 13:40:02.997862"
$ws.Range("F29").Value = "2021-10-05 13:40:02.997864"
$ws.Range("F30").Value = "2021-10-05 13:40:02.997867"
$ws.Range("F31").Value = "2021-10-05 13:40:02.997869"
$ws.Range("F32").Value = "2021-10-05 13:40:02.997872"
$ws.Range("F33").Value = "2021-10-05 13:40:02.997874"
$ws.Range("F34").Value = "2021-10-05 13:40:02.997877"
$ws.Range("F35").Value = "2021-10-05 13:40:02.997880"
$ws.Range("F36").Value = "2021-10-05 13:40:02.997882"
$ws.Range("F37").Value = "2021-10-05 13:40:02.997885"
$ws.Range("F38").Value = "2021-10-05 13:40:02.997887"
$ws.Range("F39").Value = "2021-10-05 13:40:02.997889"
$ws.Range("F40").Value = "2021-10-05 13:40:02.997892"
$ws.Range("F41").Value = "2021-10-05 13:40:02.997894"
$ws.Range("F42").Value = "2021-10-05 13:40:02.997897"
$ws.Range("F43").Value = "2021-10-05 13:40:02.997900"
$ws.Range("F44").Value = "2021-10-05 13:40:02.997902"
$ws.Range("F45").Value = "2021-10-05 13:40:02.997905"
$ws.Range("F46").Value = "2021-10-05 13:40:02.997907"
$ws.Range("F47").Value = "2021-10-05 13:40:02.997910"
$ws.Range("F48").Value = "2021-10-05 13:40:02.997912"
$ws.Range("F49").Value = "2021-10-05 13:40:02.997915"
$ws.Range("F50").Value = "2021-10-05 13:40:02.997917"
$ws.Range("F51").Value = "2021-10-05 13:40:02.997920"
$ws.Range("F52").Value = "2021-10-05 13:40:02.997922"
$ws.Range("F53").Value = "2021-10-05 13:40:02.997925"
$ws.Range("F54").Value = "2021-10-05 13:40:02.997928"
$ws.Range("F55").Value = "2021-10-05 13:40:02.997930"
$ws.Range("F56").Value = "2021-10-05 13:40:02.997933"
$ws.Range("F57").Value = "2021-10-05 13:40:02.997935"
$ws.Range("F58").Value = "2021-10-05 13:40:02.997938"
$ws.Range("F59").Value = "2021-10-05 13:40:02.997940"
$ws.Range("F60").Value = "2021-10-05 13:40:02.997943"
$ws.Range("F61").Value = "2021-10-05 13:40:02.997945"
$ws.Range("F62").Value = "2021-10-05 13:40:02.997948"
$ws.Range("F63").Value = "2021-10-05 13:40:02.997950"
$ws.Range("F64").Value = "2021-10-05 13:40:02.997953"
$ws.Range("F65").Value = "2021-10-05 13:40:02.997955"
$ws.Range("F66").Value = "2021-10-05 13:40:02.997959"
$ws.Range("F67").Value = "2021-10-05 13:40:02.997961"
$ws.Range("F68").Value = "2021-10-05 13:40:02.997964"
$ws.Range("F69").Value = "2021-10-05 13:40:02.997967"
$ws.Range("F70").Value = "2021-10-05 13:40:02.997969"
$ws.Range("F71").Value = "2021-10-05 13:40:02.997972"
$ws.Range("F72").Value = "2021-10-05 13:40:02.997974"
$ws.Range("F73").Value = "2021-10-05 13:40:02.997977"
$ws.Range("F74").Value = "2021-10-05 13:40:02.997979"
$ws.Range("F75").Value = "2021-10-05 13:40:02.997982"
$ws.Range("F76").Value = "2021-10-05 13:40:02.997985"
$ws.Range("F77").Value = "2021-10-05 13:40:02.997987"
$ws.Range("F78").Value = "2021-10-05 13:40:02.997992"
$ws.Range("F79").Value = "2021-10-05 13:40:02.997995"
$ws.Range("F80").Value = "2021-10-05 13:40:02.997997"
$ws.Range("F81").Value = "2021-10-05 13:40:02.998000"
$ws.Range("F82").Value = "2021-10-05 13:40:02.998002"
$ws.Range("F83").Value = "2021-10-05 13:40:02.998005"
$ws.Range("F84").Value = "2021-10-05 13:40:02.998008"
$ws.Range("F85").Value = "2021-10-05 13:40:02.998010"
$ws.Range("F86").Value = "2021-10-05 13:40:02.998013"
$ws.Range("F87").Value = "2021-10-05 13:40:02.998015"
$ws.Range("F88").Value = "2021-10-05 13:40:02.998018"
$ws.Range("F89").Value = "2021-10-05 13:40:02.998020"
$ws.Range("F90").Value = "2021-10-05 13:40:02.998023"
$ws.Range("F91").Value = "2021-10-05 13:40:02.998025"
$ws.Range("F92").Value = "2021-10-05 13:40:02.998028"
$ws.Range("F93").Value = "2021-10-05 13:40:02.998030"
$ws.Range("F94").Value = "2021-10-05 13:40:02.998034"
$ws.Range("F95").Value = "2021-10-05 13:40:02.998037"
$ws.Range("F96").Value = "2021-10-05 13:40:02.998040"
$ws.Range("F97").Value = "2021-10-05 13:40:02.998042"
$ws.Range("F98").Value = "2021-10-05 13:40:02.998045"

Write-Output "applied time_taken column"
